# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# The existing sheet already has headers through column AC ("Unnamed: 28");
# we append three new columns: AD=Wins, AE=Losses, AF=Ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, formatted like the rest of row 1 (bold, centered, bordered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from the neighboring header cell (AC1) onto the new
# headers so they match the existing look (bold font, center/top alignment,
# thin border) without creating new/duplicate style entries.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Every player row (2-47) gets the team's 1998 season record: 83-79-1.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 83  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 79  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 1   # AF - Ties
}
